# "Generate Report for Handoff"
#
# Status moved from "In Translation" to "Ready for handoff" and the
# handoff-related timestamps were bumped. This touches:
#   - Overview!E2,F2 (zh-cn / de-de status) + Overview!G2 (Latest HO Xliff
#     Generate Date)
#   - zh-cn!C2 (Status) + zh-cn!H2 (Latest Handoff Datetime)
#   - de-de!C2 (Status) + de-de!H2 (Latest Handoff Datetime, shares the
#     same new timestamp as the Overview's Latest HO Xliff Generate Date)
# Widening of the now-longer "Ready for handoff" status column follows.

$wb = $excel.ActiveWorkbook

$newStatus      = "Ready for handoff"
$newHoDate      = "2016-08-20 04:43:20"   # Overview!G2 & de-de!H2
$newZhHandoff   = "2016-08-20 04:43:16"   # zh-cn!H2

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHoDate

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $newZhHandoff

# --- de-de sheet ------------------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $newHoDate

# --- Column widths: the Status columns grew to fit "Ready for handoff"
$wsOverview.Range("E1:F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
